# Updated cryptos list values (price + 1h volume change) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.066.85"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.985.33"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.73"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  +0.54%  "
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +2.02%  "
$origStyle_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.09"
$ws.Range("D7").Style = $origStyle_D7
$ws.Range("E7").Value = "  +3.65%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.26%  "
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0802"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  -1.15%  "
$origStyle_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = $origStyle_D11
$ws.Range("E11").Value = "  +0.15%  "
$origStyle_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.01"
$ws.Range("D12").Style = $origStyle_D12
$ws.Range("E12").Value = "  +9.48%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$origStyle_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.846"
$ws.Range("D13").Style = $origStyle_D13
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$origStyle_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.18"
$ws.Range("D14").Style = $origStyle_D14
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "2.277.14"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("E16").Value = "  +4.10%  "
$ws.Range("D17").Value = "1.986.09"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "36.931.16"
$ws.Range("E18").Value = "  +1.26%  "
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.26"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  +2.28%  "
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.27"
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("E25").Value = "  +0.27%  "
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.153"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  +9.09%  "
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.28"
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = "  +0.88%  "
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.70"
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = "  +2.24%  "
$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.58"
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("E30").Value = "  +17.80%  "
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("E33").Value = "  +0.45%  "
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.54"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  +6.14%  "
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +0.33%  "
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.55"
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = "  -6.79%  "
$origStyle_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0985"
$ws.Range("D40").Style = $origStyle_D40
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("E43").Value = "  +1.00%  "
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.49"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  +3.10%  "
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.38"
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("D46").Value = "1.372.08"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +1.74%  "
$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.50"
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = "  +6.41%  "
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("E51").Value = "  +10.49%  "
